# code thêm tạo report lương tổng hợp
#
# 1) "Đơn sale chính" sheet: add a new HD-LUXURY order (row 4) and push the
#    existing "Tổng" (Total) row down to row 5, recalculated with the new
#    totals.
# 2) "Lương" sheet: refresh the summary salary figures that depend on the
#    new order (công, lương cơ bản, chiết khấu, tổng lương...).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Đơn sale chính" -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a brand-new row for the order so the old "Tổng" row slides from 4 -> 5
$ws1.Rows.Item(4).Insert()

$ws1.Range("A4").Value = "HD-LUXURY"
$ws1.Range("B4").Value = 641
$ws1.Range("C4").NumberFormat = "@"
$ws1.Range("C4").Value = "08-10-2024"
$ws1.Range("D4").Value = "LONG XUYÊN"
$ws1.Range("E4").Value = "Thanh nhã"
$ws1.Range("F4").Value = "Cá nhân"
$ws1.Range("G4").Value = "Tiêm botox"
$ws1.Range("H4").Value = 5000000
$ws1.Range("K4").Value = 5000000
$ws1.Range("L4").Value = 5000000
$ws1.Range("M4").Value = 0.13
$ws1.Range("N4").Value = 650000

# Update the "Tổng" (Total) row, now on row 5, with the new aggregate values
$ws1.Range("B5").Value = 3
$ws1.Range("H5").Value = 17000000
$ws1.Range("J5").Value = 0
$ws1.Range("K5").Value = 17000000
$ws1.Range("L5").Value = 16000000
$ws1.Range("M5").Value = 0
$ws1.Range("N5").Value = 1750000

# --- Sheet 2: "Lương" ------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B12").Value = 10.5
$ws2.Range("B13").Value = 3000000
$ws2.Range("B14").Value = 1750000
$ws2.Range("B32").Value = -1250000
$ws2.Range("B34").Value = -1250000
